$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest six years (2004-2009): rows 2-7. This shifts the existing
# 2010-2020 rows (formerly rows 8-18) up to rows 2-12, preserving their
# values/styles untouched.
$ws.Range("A2:A7").EntireRow.Delete() | Out-Null

# Give the new row 13 the same look (style) as the row above it (2020),
# then fill in the 2021 figures.
$ws.Range("A12").Copy($ws.Range("A13")) | Out-Null

$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = [double]"22639.22"
$ws.Cells.Item(13, 3).Value = [double]"43054.84"
$ws.Cells.Item(13, 4).Value = [double]"12248.02"
$ws.Cells.Item(13, 5).Value = [double]"20627.14"
$ws.Cells.Item(13, 6).Value = [double]"26178.59"
$ws.Cells.Item(13, 7).Value = [double]"232280.06"
